# HM_Bank_Schema.docx edits:
#   1. account_type column type:      VARCHAR(10)           -> ENUM('savings','current','zero_balance')
#   2. transaction_type label case:   transaction_type       -> Transaction_type
#   3. transaction_type column type:  VARCHAR(20) NOT NULL,  -> ENUM('deposit','withdrawal','transfer\u2019) NOT NULL,
#
# We locate each target run with Find.Execute (no in-line replacement text) and then
# assign Range.Text directly so straight apostrophes are not auto-converted into a
# matched pair of curly "smart quotes" by the find/replace auto-correct pipeline.
# The closing apostrophe before the final ")" is deliberately the Unicode right
# single quotation mark (U+2019), matching the source edit exactly.

$d = $word.ActiveDocument

# --- 1) account_type VARCHAR(10) -> ENUM('savings','current','zero_balance') ---
$r1 = $d.Content
$found1 = $r1.Find.Execute(" VARCHAR(10),", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if ($found1) {
    $r1.Text = " ENUM('savings','current','zero_balance'),"
    Write-Output "change1: ok"
} else {
    Write-Output "change1: NOT FOUND"
}

# --- 2) "    transaction_type " -> "    Transaction_type " ---
$r2 = $d.Content
$found2 = $r2.Find.Execute("    transaction_type ", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if ($found2) {
    $r2.Text = "    Transaction_type "
    Write-Output "change2: ok"
} else {
    Write-Output "change2: NOT FOUND"
}

# --- 3) transaction_type column: VARCHAR(20) NOT NULL, -> ENUM('deposit','withdrawal','transfer'+U2019+') NOT NULL, ---
$r3 = $d.Content
$found3 = $r3.Find.Execute("VARCHAR(20) NOT NULL,", $true, $false, $false, $false, $false,
                            $true, 1, $false, "", 0)
if ($found3) {
    $rightQuote = [char]0x2019
    $r3.Text = "ENUM('deposit','withdrawal','transfer" + $rightQuote + ") NOT NULL,"
    Write-Output "change3: ok"
} else {
    Write-Output "change3: NOT FOUND"
}
